# "add property of skill": insert a new "AutoAtkDis" column right after the
# existing "AtkDis" column (H), shifting NeedTar/DefaultHitTime/ShowName one
# column to the right, and fill it with 1 for every skill row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AtkDis lives in column H; NeedTar (the next header) is in column I.
# Inserting a whole column at I pushes NeedTar/DefaultHitTime/ShowName right
# and leaves a blank column I for the new field, matching how Excel's
# "Insert Column" behaves when done from the UI.
$ws.Columns("I").Insert()

# New header cell.
$ws.Range("I1").Value = "AutoAtkDis"

# All 8 skill rows (2-9) get the same default value of 1 for the new field.
$ws.Range("I2:I9").Value = 1

# Mirror the author's final selection on the newly added column.
[void]$ws.Range("I2:I9").Select()
